$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.838.68"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.230.10"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.99%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.01"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.629"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.97%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "76.31"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +5.83%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.620"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.87"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +6.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0944"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.14"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.76%  "
$ws.Range("E13").Value = "  -1.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.565.75"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.73"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.855"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.79%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.229.81"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.98%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.777.07"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0973"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.14"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.15"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.21"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "230.28"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.84%  "
$ws.Range("E24").Value = "  +0.21%  "
$ws.Range("E25").Value = "  -5.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.10"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.91%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.31"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -5.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.41"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +16.22%  "
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "169.11"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.47"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.64%  "
$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.47"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +7.73%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0848"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +4.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.120"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -4.64%  "
$ws.Range("E35").Value = "  -0.26%  "
$ws.Range("E36").Value = "  -3.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.85"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0297"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.57"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.85"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "112.11"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +14.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.203"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -4.76%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "60.11"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.79"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.94%  "
$ws.Range("E46").Value = "  -3.37%  "
$ws.Range("E47").Value = "  -0.44%  "
$ws.Range("E48").Value = "  -3.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.17"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.27"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.19"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -14.16%  "
